$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "95.336.39"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.93%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.607.38"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.83%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "2.67"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +39.54%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.00"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "222.92"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -5.89%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "636.58"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -3.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.417"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.63%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.17"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +9.73%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.00"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.602.08"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.89%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.49"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +6.37%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.211"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.46%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -9.08%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.44"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -7.20%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.281.27"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.82%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "95.330.16"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.43%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "21.95"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +17.38%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.65%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.81"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +6.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.608.50"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.68%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +7.73%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.279"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +47.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "512.27"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.52%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -7.23%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "120.88"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +14.12%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -10.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.78"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.69%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.791.41"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.67"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -6.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.75"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.65%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.41%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.08%  "

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.11%  "

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "PolygonEcosystemToken"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.614"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.98%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "32.48"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.10%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -7.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.75"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.78%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.30"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "578.28"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -9.50%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.01"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.31%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.61"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.46%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.494"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.79%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0505"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +11.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.155"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.41%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.951"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.77%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.93"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "228.79"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +11.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.84"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.18%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.51"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.52%  "
